$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Get Off My Nerves Chiropractic"
$ws.Range("B2").Value = "Capital City Chiropractic"
$ws.Range("C2").Value = "Leist Chiropractic Clinic"
$ws.Range("D2").Value = "Zachary Chiropractic Clinic"
$ws.Range("E2").Value = "Sonnier Chiropractic Clinic"
$ws.Range("F2").Value = "Louisiana Chiropractic"
$ws.Range("G2").Value = "Prewitt Chiropractic Clinic"
$ws.Range("H2").Value = "Ozark Chiropractic Clinic"
$ws.Range("I2").Value = "Spencer’s Chiropractic Clinic"
$ws.Range("J2").Value = "Community Chiropractic"
$ws.Range("K2").Value = "Chiro-Practical"
$ws.Range("L2").Value = "Family Chiropractic Clinic"
$ws.Range("M2").Value = "Massage Emporium"
$ws.Range("N2").Value = "Core Chiropractic Clinic"
$ws.Range("O2").Value = "Underwood Chiropractic Clinic"
$ws.Range("P2").Value = "Family and Sports Chiropractic"
$ws.Range("Q2").Value = "Spine and Sport Rehab Center"
$ws.Range("R2").Value = "Chiropractic Sports & Injury Center"
$ws.Range("S2").Value = "Family Chiropractic Clinic"
$ws.Range("T2").Value = "Capitol Spine and Rehabilitation"

$wb.Save()
